# Update the "Ver Vigencia de Programas" row of the test-tracking table:
# mark the use case as tested/closed (green shading + updated statuses).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row whose first cell is "Ver Vigencia de Programas".
# (Cell text carries a trailing cell-mark, so use Contains rather than -eq.)
$targetRow = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $firstCellText = $t.Cell($r, 1).Range.Text
    if ($firstCellText.Contains("Ver Vigencia de Programas")) {
        $targetRow = $r
        break
    }
}

# New text for columns 2..11 of that row.
$newValues = @(
    "Realizada",
    "N/A",
    "26 y 27/07/2020",
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "N/A",
    "Aprobado (CU cerrado - No se debe modificar)"
)

# Green fill used to highlight the now-completed row (RGB 00B050 -> BGR long).
$greenFill = 5287936

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $col = $i + 2
    $cell = $t.Cell($targetRow, $col)
    $rng = $cell.Range
    [void]$rng.MoveEnd(1, -1)
    $rng.Text = $newValues[$i]
    $cell.Shading.BackgroundPatternColor = $greenFill
}
